$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update team name "BORGES ITAQUI F.C." -> "Grêmio imortal 36"
# It appears in column E (Mandante_Nome) on row 9, 57, 72
# and column G (Visitante_Nome) on row 25, 40, 88
$ws.Range("E9").Value = "Grêmio imortal 36"
$ws.Range("G25").Value = "Grêmio imortal 36"
$ws.Range("G40").Value = "Grêmio imortal 36"
$ws.Range("E57").Value = "Grêmio imortal 36"
$ws.Range("E72").Value = "Grêmio imortal 36"
$ws.Range("G88").Value = "Grêmio imortal 36"

# Update the corresponding ID_Time values from 3914981 to 24856400
$ws.Range("F9").Value = 24856400
$ws.Range("H25").Value = 24856400
$ws.Range("H40").Value = 24856400
$ws.Range("F57").Value = 24856400
$ws.Range("F72").Value = 24856400
$ws.Range("H88").Value = 24856400
